$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the minutes worked on 07.01.2024 (row 21, column B)
$ws.Range("B21").Value = 510

# Replace the placeholder note with the final worked time ranges
$ws.Range("D21").Value = "10:30-12:00; 13:30-15:00;16:00:19:00;20:00-22:30"

# Move the active selection to B22, as after finishing data entry on row 21
$ws.Range("B22").Select()
